$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-26 Thursday", "2026-02-27 Friday"),
    @("935÷2=", "628÷9="),
    @("682÷9=", "926÷4="),
    @("393÷6=", "908÷8="),
    @("394÷7=", "977÷7="),
    @("855÷8=", "752÷6="),
    @("880÷5=", "104÷8="),
    @("492÷2=", "938÷5="),
    @("233÷6=", "650÷2="),
    @("641÷2=", "222÷4="),
    @("419÷8=", "656÷9="),
    @("328÷7=", "960÷4="),
    @("584÷7=", "692÷9="),
    @("809÷3=", "839÷9="),
    @("952÷7=", "369÷9="),
    @("502÷9=", "984÷8="),
    @("743÷6=", "150÷4="),
    @("707÷9=", "953÷3="),
    @("228÷3=", "596÷2="),
    @("233÷5=", "635÷7="),
    @("491÷7=", "671÷7="),
    @("557÷5=", "861÷2="),
    @("936÷7=", "768÷9="),
    @("551÷8=", "105÷6="),
    @("425÷7=", "291÷9="),
    @("500÷8=", "222÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
